$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "69.457.32"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  +2.69%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.395.71"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  +2.36%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  +0.03%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "589.47"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +2.23%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "181.80"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +4.60%  "

# Row 7
$ws.Range("E7").Value = "  -0.09%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.597"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  +1.48%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.195"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +9.80%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.589"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +2.35%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "48.82"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +6.82%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0000282"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +4.79%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "693.07"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  -0.65%  "

# Row 14
$ws.Range("B14").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C14").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "3.943.18"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +2.06%  "

# Row 15
$ws.Range("B15").Value = "Polkadot"
$ws.Range("C15").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "8.59"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +2.71%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "69.627.73"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +2.91%  "

# Row 17
$ws.Range("B17").Value = "WrappedEther"
$ws.Range("C17").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.415.86"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +2.97%  "

# Row 18
$ws.Range("B18").Value = "TRON"
$ws.Range("C18").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.120"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +1.83%  "

# Row 19
$ws.Range("E19").Value = "  +2.24%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.43"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +4.57%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.905"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +2.01%  "

# Row 22
$ws.Range("E22").Value = "  +1.51%  "

# Row 23
$ws.Range("E23").Value = "  +1.79%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "104.51"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +6.85%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.97"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +3.11%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.74"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +2.51%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.64"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +3.54%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "34.32"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +3.81%  "

# Row 29
$ws.Range("E29").Value = "  +2.92%  "

# Row 30
$ws.Range("E30").Value = "  +0.29%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.69"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +11.92%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "11.20"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +2.49%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "557.64"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -1.64%  "

# Row 34
$ws.Range("E34").Value = "  +2.13%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "58.20"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +1.19%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.723.35"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +0.25%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.00"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +0.06%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.142"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +9.05%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "35.02"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +2.65%  "

# Row 40
$ws.Range("E40").Value = "  +6.50%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.22"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +2.10%  "

# Row 42
$ws.Range("E42").Value = "  +2.66%  "

# Row 43
$ws.Range("E43").Value = "  +2.42%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0420"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +4.06%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.25"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -1.84%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.66"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -1.25%  "

# Row 47
$ws.Range("E47").Value = "  +1.83%  "

# Row 48
$ws.Range("B48").Value = "Mantle"
$ws.Range("C48").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.40"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +6.68%  "

# Row 49
$ws.Range("B49").Value = "FirstDigitalUSD"
$ws.Range("C49").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.00"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -0.21%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "132.80"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +3.60%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.61"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -2.36%  "
